$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-label the header row.
#    Old headers: B1=Itemname, C1=Desc, D1=Price, E1=ValidUntil
#    New headers: B1=ValidUntil, C1=Item_name, D1=Price
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "ValidUntil"
$ws.Range("C1").Value = "Item_name"
$ws.Range("D1").Value = "Price"

# Column E ("ValidUntil" in the old layout) is no longer used, drop it
# entirely - this also shrinks the used range back down to column D.
$ws.Columns.Item(5).Delete()

# ---------------------------------------------------------------------------
# 2. Update the existing data row (row 2): the item name changes.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Waste bin"

# ---------------------------------------------------------------------------
# 3. Append the rest of the catalogue rows (3-19).
#    Columns: A = running index (number), B = validity text,
#             C = item name (text), D = price (kept as text, like row 2).
# ---------------------------------------------------------------------------
$rows = @(
    @(1,  "Price valid until 28 Mar", "Scented tealight",            "1"),
    @(2,  "Price valid until 28 Mar", "18-piece cutlery set",        "1"),
    @(3,  "Price valid until 28 Mar", "Cushion",                     "1"),
    @(4,  "Price valid until 28 Mar", "Storage box with lid",        "1"),
    @(5,  "Price valid until 28 Mar", "Shower curtain",              "1"),
    @(6,  "Price valid until 28 Mar", "Hook",                        "1"),
    @(7,  "Price valid until 28 Mar", "Storage box with lid",        "3"),
    @(8,  "Price valid until 28 Mar", "Cushion cover",               "2"),
    @(9,  "Price valid until 28 Mar", "White wine glass",            "1"),
    @(10, "Price valid until 28 Mar", "Rug, low pile",               "9"),
    @(11, "Price valid until 28 Mar", "Scented candle in glass",     "1"),
    @(12, "Price valid until 28 Mar", "Rug, low pile",               "39"),
    @(13, "Price valid until 28 Mar", "Hook",                        "1"),
    @(14, "Price valid until 28 Mar", "Mirror",                      "17"),
    @(15, "Price valid until 28 Mar", "Throw",                       "12"),
    @(16, "Price valid until 28 Mar", "Throw",                       "9"),
    @(17, "Price valid until 28 Mar", "Block-out curtains, 1 pair",  "9")
)

$lastRow = 2 + $rows.Count

# Fill in the raw values first.
$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Copy row 2's formatting for column A (bold / bordered / centered style)
# onto the freshly added index cells so they match the look of the header.
$ws.Range("A2").Copy()
$ws.Range("A3:A$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select() | Out-Null
